$d = $word.ActiveDocument

# Append first new paragraph after the current last paragraph, copying the
# existing paragraph/run formatting (Word inherits it automatically via
# InsertParagraphAfter).
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Last
$newPara1.Range.InsertAfter("ELEGIMOS NO PONER UN CAMPO PRIORIDAD A LA TABLA VISIBILIDAD, DECIDIMOS QUE ORDENAMOS POR EL PRECIO DE LA MISMA, CUANTO MAS CARO ES, MAS ARRIBA VA A APARECER")

# Append second new paragraph after the first new one.
$newPara1b = $d.Paragraphs.Last
$newPara1b.Range.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Last
$newPara2.Range.InsertAfter("NO HACEMOS NINGUNA ELIMINACION FISICA DE LOS REGISTROS (DELETE QUERIES) SINO QUE MARCAMOS COMO ACTIVO FALSE EL REGISTRO EN LA TABLA CORRESPONDIENTE")

Write-Output "Inserted two paragraphs"
